# Add A/B m7 chords to the Chords sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chords")
$ws.Activate()

# Insert 4 rows before row 12 (right after the "A 7" group, before "B Maj")
$ws.Rows.Item(12).Resize(4).Insert() | Out-Null

$aM7 = @(
    @("A", "m7", "1A", 1),
    @("A", "m7", "2C", 2),
    @("A", "m7", "2E", 3),
    @("A", "m7", "2G", 5)
)
for ($i = 0; $i -lt $aM7.Length; $i++) {
    $r = 12 + $i
    $row = $aM7[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Insert 4 rows before row 26 (right after the "B 7" group, before "C Maj")
$ws.Rows.Item(26).Resize(4).Insert() | Out-Null

$bM7 = @(
    @("B", "m7", "1B", 1),
    @("B", "m7", "2D", 2),
    @("B", "m7", "2F#", 3),
    @("B", "m7", "2A", 5)
)
for ($i = 0; $i -lt $bM7.Length; $i++) {
    $r = 26 + $i
    $row = $bM7[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Reset selection to match final state
$ws.Range("D29").Select() | Out-Null
